# Applies the GitHub Actions "cryptos list" update for Mon Nov 13 13:36:07 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.909.99"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.052.13"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'245.18"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("E6").Value = "  -1.83%  "
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "'57.36"
$ws.Range("E7").Value = "  -3.08%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'59.11"
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("E10").Value = "  -3.82%  "
$ws.Range("D11").Value = "'0.0777"
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("E12").Value = "  +2.10%  "
$ws.Range("D13").Value = "'15.09"
$ws.Range("E13").Value = "  -4.70%  "
$ws.Range("D14").Value = "'0.872"
$ws.Range("E14").Value = "  +3.86%  "
$ws.Range("D15").Value = "2.351.13"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "'5.57"
$ws.Range("E16").Value = "  -3.17%  "
$ws.Range("D17").Value = "2.081.31"
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").Value = "36.841.50"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "'17.33"
$ws.Range("E19").Value = "  -3.86%  "
$ws.Range("D20").Value = "'73.09"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("D21").Value = "0.0₃0887"
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("D22").Value = "'5.40"
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("D23").Value = "'235.87"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D25").Value = "'2.45"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "'10.08"
$ws.Range("E26").Value = "  +7.44%  "
$ws.Range("D27").Value = "'2.19"
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("D28").Value = "'168.70"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").Value = "'20.02"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "'5.46"
$ws.Range("E30").Value = "  +14.85%  "
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("E32").Value = "  -1.48%  "
$ws.Range("E33").Value = "  +6.38%  "
$ws.Range("D34").Value = "'0.0612"
$ws.Range("E34").Value = "  -1.94%  "
$ws.Range("E35").Value = "  +5.50%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'1.85"
$ws.Range("E37").Value = "  +5.41%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.0849"
$ws.Range("E38").Value = "  -5.39%  "
$ws.Range("D39").Value = "'1.31"
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("E41").Value = "  -6.94%  "
$ws.Range("D42").Value = "'4.86"
$ws.Range("E42").Value = "  -6.45%  "
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("D44").Value = "'0.0955"
$ws.Range("E44").Value = "  -10.20%  "
$ws.Range("D45").Value = "'96.87"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").Value = "'16.62"
$ws.Range("E46").Value = "  -4.49%  "
$ws.Range("D47").Value = "1.304.87"
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("E48").Value = "  -4.25%  "
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("D50").Value = "'6.73"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").Value = "2.236.55"
$ws.Range("E51").Value = "  +0.16%  "
